$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 1.92
$ws.Cells.Item(2, 7).Value = 1.94
$ws.Cells.Item(2, 8).Value = 4.5
$ws.Cells.Item(2, 9).Value = 4.7
$ws.Cells.Item(2, 11).Value = 3.75
$ws.Cells.Item(2, 12).Value = 1.47
$ws.Cells.Item(2, 14).Value = 3.5
$ws.Cells.Item(2, 16).Value = 1.83
$ws.Cells.Item(2, 17).Value = 2.16
$ws.Cells.Item(2, 19).Value = 4
$ws.Cells.Item(2, 20).Value = 1.97
$ws.Cells.Item(2, 21).Value = 1.98
$ws.Cells.Item(2, 22).Value = 1.27
$ws.Cells.Item(2, 23).Value = 2.04
$ws.Cells.Item(2, 24).Value = 13
$ws.Cells.Item(2, 25).Value = 15
$ws.Cells.Item(2, 26).Value = 34
$ws.Cells.Item(2, 27).Value = 110
$ws.Cells.Item(2, 29).Value = 8
$ws.Cells.Item(2, 30).Value = 18
$ws.Cells.Item(2, 31).Value = 65
$ws.Cells.Item(2, 33).Value = 10
$ws.Cells.Item(2, 39).Value = 130
$ws.Cells.Item(2, 40).Value = 16.5
$ws.Cells.Item(3, 6).Value = 1.36
$ws.Cells.Item(3, 7).Value = 1.48
$ws.Cells.Item(3, 8).Value = 6.4
$ws.Cells.Item(3, 9).Value = 11
$ws.Cells.Item(3, 10).Value = 4.5
$ws.Cells.Item(3, 11).Value = 7.6
$ws.Cells.Item(3, 12).Value = 1.24
$ws.Cells.Item(3, 14).Value = 4.8
$ws.Cells.Item(3, 15).Value = 1.15
$ws.Cells.Item(3, 16).Value = 2.62
$ws.Cells.Item(3, 17).Value = 1.46
$ws.Cells.Item(3, 18).Value = 1.66
$ws.Cells.Item(3, 19).Value = 2.16
$ws.Cells.Item(3, 20).Value = 1.72
$ws.Cells.Item(3, 21).Value = 2.06
$ws.Cells.Item(3, 22).Value = 1.11
$ws.Cells.Item(3, 23).Value = 3.05
$ws.Cells.Item(3, 28).Value = 1000
$ws.Cells.Item(3, 29).Value = 42
$ws.Cells.Item(3, 33).Value = 1000
$ws.Cells.Item(3, 36).Value = 1000
$ws.Cells.Item(3, 40).Value = 15
$ws.Cells.Item(4, 6).Value = 1.9
$ws.Cells.Item(4, 7).Value = 2.02
$ws.Cells.Item(4, 8).Value = 4.6
$ws.Cells.Item(4, 9).Value = 5.5
$ws.Cells.Item(4, 10).Value = 3.15
$ws.Cells.Item(4, 11).Value = 3.6
$ws.Cells.Item(4, 12).Value = 1.54
$ws.Cells.Item(4, 15).Value = 1.44
$ws.Cells.Item(4, 16).Value = 1.64
$ws.Cells.Item(4, 17).Value = 2.36
$ws.Cells.Item(4, 19).Value = 4.5
$ws.Cells.Item(4, 20).Value = 2
$ws.Cells.Item(4, 21).Value = 1.8
$ws.Cells.Item(4, 22).Value = 1.22
$ws.Cells.Item(4, 23).Value = 1.98
$ws.Cells.Item(4, 25).Value = 14.5
$ws.Cells.Item(4, 26).Value = 120
$ws.Cells.Item(4, 28).Value = 7.6
$ws.Cells.Item(4, 29).Value = 7.8
$ws.Cells.Item(4, 30).Value = 21
$ws.Cells.Item(4, 33).Value = 11
$ws.Cells.Item(4, 34).Value = 24
$ws.Cells.Item(4, 36).Value = 25
$ws.Cells.Item(4, 37).Value = 26
$ws.Cells.Item(4, 38).Value = 150
$ws.Cells.Item(4, 40).Value = 22
$ws.Cells.Item(5, 6).Value = 2.86
$ws.Cells.Item(5, 7).Value = 3.2
$ws.Cells.Item(5, 8).Value = 2.58
$ws.Cells.Item(5, 9).Value = 2.86
$ws.Cells.Item(5, 10).Value = 3.1
$ws.Cells.Item(5, 11).Value = 3.5
$ws.Cells.Item(5, 12).Value = 1.47
$ws.Cells.Item(5, 14).Value = 3.2
$ws.Cells.Item(5, 15).Value = 1.39
$ws.Cells.Item(5, 16).Value = 1.73
$ws.Cells.Item(5, 17).Value = 2.2
$ws.Cells.Item(5, 19).Value = 4.2
$ws.Cells.Item(5, 20).Value = 1.81
$ws.Cells.Item(5, 21).Value = 2
$ws.Cells.Item(5, 22).Value = 1.53
$ws.Cells.Item(5, 23).Value = 1.47
$ws.Cells.Item(5, 25).Value = 1000
$ws.Cells.Item(5, 27).Value = 900
$ws.Cells.Item(5, 28).Value = 1000
$ws.Cells.Item(5, 29).Value = 14
$ws.Cells.Item(5, 30).Value = 25
$ws.Cells.Item(5, 33).Value = 1000
$ws.Cells.Item(5, 34).Value = 990
$ws.Cells.Item(5, 36).Value = 1000
$ws.Cells.Item(6, 7).Value = 1.85
$ws.Cells.Item(6, 8).Value = 4.5
$ws.Cells.Item(6, 9).Value = 5.7
$ws.Cells.Item(6, 10).Value = 3.9
$ws.Cells.Item(6, 11).Value = 4.9
$ws.Cells.Item(6, 12).Value = 1.31
$ws.Cells.Item(6, 14).Value = 4.5
$ws.Cells.Item(6, 15).Value = 1.22
$ws.Cells.Item(6, 16).Value = 2.26
$ws.Cells.Item(6, 17).Value = 1.65
$ws.Cells.Item(6, 18).Value = 1.48
$ws.Cells.Item(6, 19).Value = 2.62
$ws.Cells.Item(6, 20).Value = 1.67
$ws.Cells.Item(6, 23).Value = 2.16
$ws.Cells.Item(6, 36).Value = 1000
$ws.Cells.Item(6, 37).Value = 1000
$ws.Cells.Item(7, 6).Value = 5.7
$ws.Cells.Item(7, 7).Value = 5.8
$ws.Cells.Item(7, 8).Value = 1.83
$ws.Cells.Item(7, 12).Value = 1.52
$ws.Cells.Item(7, 13).Value = 1.1
$ws.Cells.Item(7, 14).Value = 3.15
$ws.Cells.Item(7, 15).Value = 1.45
$ws.Cells.Item(7, 16).Value = 1.71
$ws.Cells.Item(7, 17).Value = 2.36
$ws.Cells.Item(7, 18).Value = 1.26
$ws.Cells.Item(7, 19).Value = 4.5
$ws.Cells.Item(7, 20).Value = 2.12
$ws.Cells.Item(7, 21).Value = 1.83
$ws.Cells.Item(7, 22).Value = 2.18
$ws.Cells.Item(7, 24).Value = 10.5
$ws.Cells.Item(7, 25).Value = 7
$ws.Cells.Item(7, 27).Value = 19
$ws.Cells.Item(7, 28).Value = 16
$ws.Cells.Item(7, 29).Value = 7.8
$ws.Cells.Item(7, 30).Value = 10.5
$ws.Cells.Item(7, 31).Value = 22
$ws.Cells.Item(7, 33).Value = 22
$ws.Cells.Item(7, 34).Value = 24
$ws.Cells.Item(7, 35).Value = 48
$ws.Cells.Item(7, 36).Value = 160
$ws.Cells.Item(7, 37).Value = 90
$ws.Cells.Item(7, 38).Value = 110
$ws.Cells.Item(7, 39).Value = 180
$ws.Cells.Item(7, 40).Value = 160
$ws.Cells.Item(7, 41).Value = 15.5
$ws.Cells.Item(8, 6).Value = 2.76
$ws.Cells.Item(8, 7).Value = 3.1
$ws.Cells.Item(8, 8).Value = 2.46
$ws.Cells.Item(8, 9).Value = 2.7
$ws.Cells.Item(8, 10).Value = 3.45
$ws.Cells.Item(8, 11).Value = 4
$ws.Cells.Item(8, 12).Value = 1.36
$ws.Cells.Item(8, 13).Value = 1.05
$ws.Cells.Item(8, 14).Value = 4.2
$ws.Cells.Item(8, 15).Value = 1.25
$ws.Cells.Item(8, 16).Value = 2.1
$ws.Cells.Item(8, 17).Value = 1.77
$ws.Cells.Item(8, 18).Value = 1.43
$ws.Cells.Item(8, 19).Value = 2.9
$ws.Cells.Item(8, 20).Value = 1.63
$ws.Cells.Item(8, 21).Value = 2.32
$ws.Cells.Item(8, 22).Value = 1.59
$ws.Cells.Item(8, 24).Value = 18
$ws.Cells.Item(8, 25).Value = 13.5
$ws.Cells.Item(8, 26).Value = 19
$ws.Cells.Item(8, 27).Value = 38
$ws.Cells.Item(8, 28).Value = 14
$ws.Cells.Item(8, 29).Value = 9
$ws.Cells.Item(8, 30).Value = 12.5
$ws.Cells.Item(8, 31).Value = 27
$ws.Cells.Item(8, 32).Value = 22
$ws.Cells.Item(8, 33).Value = 13.5
$ws.Cells.Item(8, 34).Value = 16.5
$ws.Cells.Item(8, 35).Value = 36
$ws.Cells.Item(8, 36).Value = 46
$ws.Cells.Item(8, 37).Value = 32
$ws.Cells.Item(8, 38).Value = 40
$ws.Cells.Item(8, 39).Value = 75
$ws.Cells.Item(8, 40).Value = 28
$ws.Cells.Item(8, 41).Value = 23
$ws.Cells.Item(9, 6).Value = 5.8
$ws.Cells.Item(9, 7).Value = 5.9
$ws.Cells.Item(9, 8).Value = 1.81
$ws.Cells.Item(9, 9).Value = 1.82
$ws.Cells.Item(9, 10).Value = 3.6
$ws.Cells.Item(9, 11).Value = 3.65
$ws.Cells.Item(9, 14).Value = 3.05
$ws.Cells.Item(9, 16).Value = 1.68
$ws.Cells.Item(9, 20).Value = 2.22
$ws.Cells.Item(9, 21).Value = 1.77
$ws.Cells.Item(9, 22).Value = 2.22
$ws.Cells.Item(9, 24).Value = 9.6
$ws.Cells.Item(9, 28).Value = 15.5
$ws.Cells.Item(9, 40).Value = 150
$ws.Cells.Item(9, 41).Value = 16.5
$ws.Cells.Item(10, 10).Value = 3.95
$ws.Cells.Item(10, 12).Value = 1.44
$ws.Cells.Item(10, 14).Value = 3.75
$ws.Cells.Item(10, 16).Value = 1.95
$ws.Cells.Item(10, 17).Value = 2.04
$ws.Cells.Item(10, 19).Value = 3.7
$ws.Cells.Item(10, 24).Value = 13.5
$ws.Cells.Item(10, 26).Value = 9.800000000000001
$ws.Cells.Item(10, 28).Value = 18
$ws.Cells.Item(10, 40).Value = 95
$ws.Cells.Item(10, 41).Value = 12
$ws.Cells.Item(11, 6).Value = 3.1
$ws.Cells.Item(11, 7).Value = 3.2
$ws.Cells.Item(11, 8).Value = 2.44
$ws.Cells.Item(11, 9).Value = 2.48
$ws.Cells.Item(11, 11).Value = 3.6
$ws.Cells.Item(11, 12).Value = 1.44
$ws.Cells.Item(11, 14).Value = 3.65
$ws.Cells.Item(11, 15).Value = 1.34
$ws.Cells.Item(11, 16).Value = 1.89
$ws.Cells.Item(11, 17).Value = 2.04
$ws.Cells.Item(11, 18).Value = 1.33
$ws.Cells.Item(11, 19).Value = 3.75
$ws.Cells.Item(11, 20).Value = 1.82
$ws.Cells.Item(11, 22).Value = 1.67
$ws.Cells.Item(11, 23).Value = 1.45
$ws.Cells.Item(11, 26).Value = 15.5
$ws.Cells.Item(11, 27).Value = 36
$ws.Cells.Item(11, 28).Value = 12.5
$ws.Cells.Item(11, 30).Value = 11.5
$ws.Cells.Item(11, 31).Value = 27
$ws.Cells.Item(11, 32).Value = 21
$ws.Cells.Item(11, 33).Value = 13.5
$ws.Cells.Item(11, 36).Value = 55
$ws.Cells.Item(11, 37).Value = 38
$ws.Cells.Item(11, 39).Value = 580
$ws.Cells.Item(11, 40).Value = 36
$ws.Cells.Item(11, 41).Value = 23
$ws.Cells.Item(12, 10).Value = 9.199999999999999
$ws.Cells.Item(12, 12).Value = 1.2
$ws.Cells.Item(12, 14).Value = 8.199999999999999
$ws.Cells.Item(12, 18).Value = 2.02
$ws.Cells.Item(12, 19).Value = 1.79
$ws.Cells.Item(12, 20).Value = 2.16
$ws.Cells.Item(12, 21).Value = 1.69
$ws.Cells.Item(12, 22).Value = 1.03
$ws.Cells.Item(12, 38).Value = 48
$ws.Cells.Item(12, 40).Value = 2.66
$ws.Cells.Item(13, 6).Value = 3.3
$ws.Cells.Item(13, 7).Value = 3.35
$ws.Cells.Item(13, 9).Value = 2.42
$ws.Cells.Item(13, 12).Value = 1.44
$ws.Cells.Item(13, 14).Value = 3.9
$ws.Cells.Item(13, 16).Value = 1.95
$ws.Cells.Item(13, 17).Value = 2.02
$ws.Cells.Item(13, 18).Value = 1.37
$ws.Cells.Item(13, 20).Value = 1.8
$ws.Cells.Item(13, 22).Value = 1.71
$ws.Cells.Item(13, 24).Value = 14
$ws.Cells.Item(13, 25).Value = 10.5
$ws.Cells.Item(13, 27).Value = 30
$ws.Cells.Item(13, 28).Value = 12.5
$ws.Cells.Item(13, 31).Value = 23
$ws.Cells.Item(13, 33).Value = 14
$ws.Cells.Item(13, 34).Value = 17
$ws.Cells.Item(13, 35).Value = 36
$ws.Cells.Item(13, 36).Value = 50
$ws.Cells.Item(13, 38).Value = 46
$ws.Cells.Item(13, 39).Value = 75
$ws.Cells.Item(13, 41).Value = 21
$ws.Cells.Item(14, 6).Value = 1.56
$ws.Cells.Item(14, 8).Value = 7.2
$ws.Cells.Item(14, 9).Value = 7.6
$ws.Cells.Item(14, 11).Value = 4.5
$ws.Cells.Item(14, 12).Value = 1.33
$ws.Cells.Item(14, 14).Value = 5.1
$ws.Cells.Item(14, 16).Value = 2.32
$ws.Cells.Item(14, 18).Value = 1.53
$ws.Cells.Item(14, 19).Value = 2.76
$ws.Cells.Item(14, 20).Value = 1.82
$ws.Cells.Item(14, 24).Value = 20
$ws.Cells.Item(14, 26).Value = 1000
$ws.Cells.Item(14, 27).Value = 1000
$ws.Cells.Item(14, 28).Value = 10
$ws.Cells.Item(14, 31).Value = 200
$ws.Cells.Item(14, 32).Value = 9.6
$ws.Cells.Item(14, 35).Value = 85
$ws.Cells.Item(14, 41).Value = 320
